# Restore D9 ("To" value for rule R20) from 17 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("D9").Value = 1
